$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

# New column D ("BPA") added to the solar sheet: header + constant 85.2 for
# every data row (2-85), matching the wind sheet's BPA column label.
$ws.Range("D1").Value = "BPA"

for ($r = 2; $r -le 85; $r++) {
    $ws.Cells.Item($r, 4).Value = 85.2
}

# Match existing center/center alignment used by the sheet's other data
# columns (same visual style as column C).
$rng = $ws.Range("D1:D85")
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# Leave the selection where the saved workbook shows it.
$ws.Range("R19").Select()
